$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in E1 to reflect the new name
$ws.Range("E1").Value = "Type ID (Select)"

# Move the active selection to F7 (matches the saved view state)
$ws.Range("F7").Select()
